$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1655.56
$ws.Range("I40").Value = 1515
$ws.Range("J40").Value = 1905.4445
$ws.Range("K40").Value = 1515
$ws.Range("L40").Value = 1905.4445
$ws.Range("M40").Value = -1340
$ws.Range("N40").Value = -2255.4445

$ws.Range("H62").Value = 2371.2856
$ws.Range("I62").Value = 5500
$ws.Range("J62").Value = 1119.8
$ws.Range("K62").Value = 5500
$ws.Range("L62").Value = 1119.8
$ws.Range("M62").Value = -4876
$ws.Range("N62").Value = -2367.8

$ws.Range("H65").Value = 2371.2856
$ws.Range("I65").Value = 5500
$ws.Range("J65").Value = 1119.8
$ws.Range("K65").Value = 27500
$ws.Range("L65").Value = 5599
$ws.Range("M65").Value = -24380
$ws.Range("N65").Value = -11839

$ws.Range("H86").Value = 3891.2307
$ws.Range("I86").Value = 3881
$ws.Range("J86").Value = 3900
$ws.Range("K86").Value = 3881
$ws.Range("L86").Value = 3900
$ws.Range("M86").Value = -2758
$ws.Range("N86").Value = -6146

$ws.Range("H89").Value = 3891.2307
$ws.Range("I89").Value = 3881
$ws.Range("J89").Value = 3900
$ws.Range("K89").Value = 19405
$ws.Range("L89").Value = 19500
$ws.Range("M89").Value = -13789
$ws.Range("N89").Value = -30732

$ws.Range("H92").Value = 1785.76
$ws.Range("I92").Value = 1602.0952
$ws.Range("K92").Value = 1602.0952
$ws.Range("M92").Value = -354.0952

$ws.Range("H99").Value = 78474.16
$ws.Range("I99").Value = 742
$ws.Range("J99").Value = 169161.67
$ws.Range("K99").Value = 2226
$ws.Range("L99").Value = 507485.01
$ws.Range("M99").Value = -728
$ws.Range("N99").Value = -510481.01

$ws.Range("H116").Value = 3279.8
$ws.Range("I116").Value = 3700
$ws.Range("K116").Value = 3700
$ws.Range("M116").Value = -258

$ws.Range("H132").Value = 3714.1292
$ws.Range("I132").Value = 3471.2666
$ws.Range("J132").Value = 11000
$ws.Range("K132").Value = 10413.7998
$ws.Range("L132").Value = 33000
$ws.Range("M132").Value = -7883.799800000001
$ws.Range("N132").Value = -38060

$ws.Range("H138").Value = 2871.9546
$ws.Range("J138").Value = 3454.5454
$ws.Range("L138").Value = 10363.6362
$ws.Range("N138").Value = -20643.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6262.1284
$ws.Range("I32").Value = 6628.394
$ws.Range("J32").Value = 4247.6665
$ws.Range("K32").Value = 6628.394
$ws.Range("L32").Value = 4247.6665
$ws.Range("M32").Value = -6341.394
$ws.Range("N32").Value = -4821.6665

$ws.Range("H45").Value = 1889
$ws.Range("I45").Value = 1862.3182
$ws.Range("J45").Value = 1962.375
$ws.Range("K45").Value = 1862.3182
$ws.Range("L45").Value = 1962.375
$ws.Range("M45").Value = -1485.3182
$ws.Range("N45").Value = -2716.375

$ws.Range("H74").Value = 14288500
$ws.Range("I74").Value = 20001300
$ws.Range("J74").Value = 6497
$ws.Range("K74").Value = 20001300
$ws.Range("L74").Value = 6497
$ws.Range("M74").Value = -20000426
$ws.Range("N74").Value = -8245

$ws.Range("H77").Value = 14288500
$ws.Range("I77").Value = 20001300
$ws.Range("J77").Value = 6497
$ws.Range("K77").Value = 100006500
$ws.Range("L77").Value = 32485
$ws.Range("M77").Value = -100002132
$ws.Range("N77").Value = -41221

$ws.Range("H97").Value = 6892.6875
$ws.Range("I97").Value = 8956
$ws.Range("J97").Value = 702.75
$ws.Range("K97").Value = 8956
$ws.Range("L97").Value = 702.75
$ws.Range("M97").Value = -8460
$ws.Range("N97").Value = -1694.75

$ws.Range("H102").Value = 1566.25
$ws.Range("I102").Value = 1208.6364
$ws.Range("J102").Value = 5500
$ws.Range("K102").Value = 1208.6364
$ws.Range("L102").Value = 5500
$ws.Range("M102").Value = 413.3635999999999
$ws.Range("N102").Value = -8744

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1072.05
$ws.Range("I20").Value = 1129.4
$ws.Range("J20").Value = 900
$ws.Range("K20").Value = 1129.4
$ws.Range("L20").Value = 900
$ws.Range("M20").Value = -882.4000000000001
$ws.Range("N20").Value = -1394

$ws.Range("H99").Value = 920.5
$ws.Range("I99").Value = 832.25
$ws.Range("J99").Value = 1450
$ws.Range("K99").Value = 832.25
$ws.Range("L99").Value = 1450
$ws.Range("M99").Value = 665.75
$ws.Range("N99").Value = -4446

$ws.Range("H107").Value = 858.4286
$ws.Range("I107").Value = 917.2105
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 917.2105
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 1002.7895
$ws.Range("N107").Value = -4140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2195.6
$ws.Range("I58").Value = 536.6429000000001
$ws.Range("J58").Value = 6066.5
$ws.Range("K58").Value = 536.6429000000001
$ws.Range("L58").Value = 6066.5
$ws.Range("M58").Value = -333.6429000000001
$ws.Range("N58").Value = -6472.5

$ws.Range("H132").Value = 2732.3125
$ws.Range("I132").Value = 1958.5385
$ws.Range("J132").Value = 6085.3335
$ws.Range("K132").Value = 5875.6155
$ws.Range("L132").Value = 18256.0005
$ws.Range("M132").Value = -3345.6155
$ws.Range("N132").Value = -23316.0005

$ws.Range("H136").Value = 2195.6
$ws.Range("I136").Value = 536.6429000000001
$ws.Range("J136").Value = 6066.5
$ws.Range("K136").Value = 1609.9287
$ws.Range("L136").Value = 18199.5
$ws.Range("M136").Value = 940.0712999999998
$ws.Range("N136").Value = -23299.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 67857
$ws.Range("J133").Value = 67857
$ws.Range("L133").Value = 67857
$ws.Range("N133").Value = -77977

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2468.3125
$ws.Range("I16").Value = 2468.3125
$ws.Range("K16").Value = 2468.3125
$ws.Range("M16").Value = -2298.3125

$ws.Range("H46").Value = 1023.4667
$ws.Range("I46").Value = 836
$ws.Range("J46").Value = 1187.5
$ws.Range("K46").Value = 836
$ws.Range("L46").Value = 1187.5
$ws.Range("M46").Value = -648
$ws.Range("N46").Value = -1563.5

$ws.Range("H93").Value = 1908.3334
$ws.Range("I93").Value = 1842.8572
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 1842.8572
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -594.8571999999999
$ws.Range("N93").Value = -4496

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1800
$ws.Range("I122").Value = 1574.862
$ws.Range("J122").Value = 2888.1667
$ws.Range("K122").Value = 4724.586
$ws.Range("L122").Value = 8664.500100000001
$ws.Range("M122").Value = -2274.586
$ws.Range("N122").Value = -13564.5001
